$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Verif fonction utilisée debloquée" - the Chap/difficulty value (column A)
# for the unlocked quest rows (2-5) drops from 4 to 3.
$ws.Range("A2:A5").Value = 3

# Window/view bookkeeping captured by the same save: the user had scrolled
# to cell A6 (instead of H7) and zoomed out to 60% (instead of 93%).
$ws.Range("A6").Select()
$excel.ActiveWindow.Zoom = 60
